# Apply crypto price/volume updates from the GitHub Actions scrape run.
# Every written cell is forced to text (apostrophe-prefix) and then restyled
# back to the workbook's default "Normal" style so numeric-looking strings
# (e.g. "20.12", "68.310.40") stay literal text instead of being coerced to
# real numbers/dates by Excel's type inference, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'68.310.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.28%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.913.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.77%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'485.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.49%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'146.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.83%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.76%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  -0.12%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.735"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.39%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.04%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.0000344"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.22%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'43.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.12%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'10.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.06%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'4.530.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.83%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'3.907.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.81%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'14.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.31%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "'  -1.19%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'20.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.32%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  -1.14%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'68.324.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.09%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'433.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.05%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'3.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.69%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'15.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.45%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'88.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.71%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'11.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +15.00%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'11.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +10.83%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'3.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.06%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'38.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.63%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'5.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.44%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("B30").Value = "'Bittensor"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'712.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.89%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("B31").Value = "'Cosmos"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'13.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.74%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  -2.36%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  +4.38%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'6.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +16.24%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'41.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.61%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("B36").Value = "'PEPE"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.0₃0873"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.73%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("B37").Value = "'OKB"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'61.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +4.68%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("B38").Value = "'Kaspa"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.146"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.98%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("B39").Value = "'Dai"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.10%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.391"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +15.68%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.0490"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.04%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'2.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +15.68%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = "'  +2.25%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'2.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.67%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'  -2.11%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "'FirstDigitalUSD"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.06%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = "'ApeXProtocol"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'3.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.16%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'3.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.85%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -4.04%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'145.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.86%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.0₆0335"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +30.93%  "
$ws.Range("E51").Style = "Normal"
